$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11395
$ws.Range("C3").Value = 6092
$ws.Range("C4").Value = 9007
$ws.Range("C5").Value = 6835
$ws.Range("C6").Value = 5243
$ws.Range("C7").Value = 8462
$ws.Range("C8").Value = 22730
$ws.Range("C9").Value = 16750
$ws.Range("C10").Value = 5036
$ws.Range("C11").Value = 3992
$ws.Range("C12").Value = 59
